# This workbook (CDFW Fish Bulletin 108, Table 34) contains OCR-garbled
# species names in its "port landings" breakdown tables. This edit corrects
# several of those garbled names to their clean equivalents (matching the
# clean spellings already used elsewhere in the sheet, e.g. in the Eureka
# port table), and updates the sheet's active selection/scroll position to
# reflect where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fort Bragg port table (rows 21-33)
$ws.Range("B22").Value = "Albacore"        # was "A Ilia core"
$ws.Range("B24").Value = "Rockfish"        # was "Rockfbh"
$ws.Range("B27").Value = "English sole "   # was "English s<i|e "

# Fields Landing port table (rows 43-55)
$ws.Range("B44").Value = "Petrale sole "   # was "I'etrale sole    1"
$ws.Range("B46").Value = "Crab"            # was "Crab   i"
$ws.Range("B47").Value = "Rockfish"        # was "Kockfish"
$ws.Range("B48").Value = "Rex sole"        # was "lb x sole "
$ws.Range("B53").Value = "OCean shrimp"    # was "(Wan shrimp "

# Update the view: scroll down to around row 39 and select B51, matching
# where work continued in the port-level landings tables.
$ws.Activate()
$ws.Range("A39").Select()
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B51").Select()
